$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose "relevant" flag (column B) flips from 1 to 0 -- fewer Wikipedia
# pages counted as used/relevant, per the commit message.
$rowsToZero = @(49, 84, 123, 167, 169, 172, 173, 186, 194, 195, 196, 197, 203, 213, 239, 278, 320, 323, 328, 353, 354, 390, 391)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Give column B an explicit width (close best-fit rendering of the "relevant"
# header / 0-1 values), matching the new <col> entry for column B.
$ws.Columns.Item(2).ColumnWidth = 11.83

# Turn on AutoFilter over the full used range -- this both adds the
# <autoFilter> element on the sheet and the workbook-level hidden
# _FilterDatabase defined name.
$ws.Range("A1:B450").AutoFilter() | Out-Null

$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$450")
$fdb.Visible = $false

# Move the selection / cursor back to the top of the sheet and drop the
# scrolled-down view (topLeftCell) that was left over from browsing.
$ws.Activate()
$ws.Range("A5").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
